$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column width updates (D: 13->12, E: 24->22, F: 24->26).
# ColumnWidth setter stores width+0.8333333333333334 in the OOXML "width"
# attribute, so subtract that offset to land on the exact target values.
$ws.Columns.Item(4).ColumnWidth = 12 - 0.8333333333333334
$ws.Columns.Item(5).ColumnWidth = 22 - 0.8333333333333334
$ws.Columns.Item(6).ColumnWidth = 26 - 0.8333333333333334

# Data rows 2-14 (GRUPO, PRESUPUESTO, VENTA, POR CUMPLIR, CUMPLIMIENTO)
$data = @(
    @("240X120 PORCELANATO", 1680.23389242503, 0, 1680.23389242503, 0),
    @("240X80 PORCELANATO", 5504.61890386263, 0, 5504.61890386263, 0),
    @("FREGADEROS DE COCINA", 1304.0286065816, 0, 1304.0286065816, 0),
    @("GRIFERIAS", 150, 0, 150, 0),
    @("INODOROS", 849.84419682004, 0, 849.84419682004, 0),
    @("LAVABOS", 709.368813030059, 0, 709.368813030059, 0),
    @("NO RESURTIBLES", 516.121873547834, 0, 516.121873547834, 0),
    @("OTROS", 0, 0, 0, 0),
    @("PANELES DECORATIVOS", 388.107983534392, 0, 388.107983534392, 0),
    @("PIEDRA SINTERIZADA", 3506.66949822329, 0, 3506.66949822329, 0),
    @("PORCELANATO", 32404.8, 253.96, 32150.84, 0.007837110551523231),
    @("PUERTAS DE SEGURIDAD", 1332.52398144409, 0, 1332.52398144409, 0),
    @("SAL SOLUBLE", 2137.44930155624, 705.02, 1432.42930155624, 0.3298417415031492)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $ws.Cells.Item($row, 5).Value = $item[3]
    $ws.Cells.Item($row, 6).Value = $item[4]
    $row++
}

# Row 15 = TOTAL row
$ws.Cells.Item(15, 2).Value = "TOTAL"
$ws.Cells.Item(15, 2).HorizontalAlignment = -4152
$ws.Cells.Item(15, 3).Value = 50483.7670510252
$ws.Cells.Item(15, 4).Value = 958.98
$ws.Cells.Item(15, 5).Value = 49524.7870510252
$ws.Cells.Item(15, 6).Value = 0.0189958090692942

# Clear cell A15 (previously held "LINDAO ZUÑIGA BRYAN JOSE" from old PIEDRA SINTERIZADA row)
$ws.Cells.Item(15, 1).ClearContents()

# Delete rows 16-19 (old PORCELANATO, PUERTAS DE SEGURIDAD, SAL SOLUBLE, TOTAL rows) -- now beyond new range
$ws.Range("A16:F19").Delete()
